$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1528.9474
$ws.Range("I107").Value = 1575
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 1575
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 345
$ws.Range("N107").Value = -4540

$ws.Range("H132").Value = 2145.4482
$ws.Range("I132").Value = 2058.2354
$ws.Range("J132").Value = 2269
$ws.Range("K132").Value = 6174.706200000001
$ws.Range("L132").Value = 6807
$ws.Range("M132").Value = -3644.706200000001
$ws.Range("N132").Value = -11867

$ws.Range("H138").Value = 2476.1313
$ws.Range("I138").Value = 1167.25
$ws.Range("J138").Value = 2728.4458
$ws.Range("K138").Value = 3501.75
$ws.Range("L138").Value = 8185.3374
$ws.Range("M138").Value = 1638.25
$ws.Range("N138").Value = -18465.3374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19231992
$ws.Range("J32").Value = 2833
$ws.Range("L32").Value = 2833
$ws.Range("N32").Value = -3407

$ws.Range("H61").Value = 3926.8
$ws.Range("I61").Value = 3343.3333
$ws.Range("K61").Value = 3343.3333
$ws.Range("M61").Value = -3131.3333

$ws.Range("H74").Value = 1562.0454
$ws.Range("I74").Value = 1471.1818
$ws.Range("J74").Value = 1834.6364
$ws.Range("K74").Value = 1471.1818
$ws.Range("L74").Value = 1834.6364
$ws.Range("M74").Value = -597.1818000000001
$ws.Range("N74").Value = -3582.6364

$ws.Range("H77").Value = 1562.0454
$ws.Range("I77").Value = 1471.1818
$ws.Range("J77").Value = 1834.6364
$ws.Range("K77").Value = 7355.909000000001
$ws.Range("L77").Value = 9173.182000000001
$ws.Range("M77").Value = -2987.909000000001
$ws.Range("N77").Value = -17909.182

$ws.Range("H97").Value = 2220.0667
$ws.Range("I97").Value = 1637.8096
$ws.Range("J97").Value = 3578.6667
$ws.Range("K97").Value = 1637.8096
$ws.Range("L97").Value = 3578.6667
$ws.Range("M97").Value = -1141.8096
$ws.Range("N97").Value = -4570.6667

$ws.Range("H122").Value = 6385.3105
$ws.Range("I122").Value = 5479.5454
$ws.Range("K122").Value = 16438.6362
$ws.Range("M122").Value = -13988.6362

$ws.Range("H126").Value = 9997
$ws.Range("I126").Value = 9997
$ws.Range("K126").Value = 29991
$ws.Range("M126").Value = -27521

$ws.Range("H129").Value = 117911.6
$ws.Range("J129").Value = 117911.6
$ws.Range("L129").Value = 117911.6
$ws.Range("N129").Value = -127911.6

$ws.Range("H132").Value = 3880.3157
$ws.Range("I132").Value = 3915.5334
$ws.Range("J132").Value = 3748.25
$ws.Range("K132").Value = 11746.6002
$ws.Range("L132").Value = 11244.75
$ws.Range("M132").Value = -9216.600199999999
$ws.Range("N132").Value = -16304.75

$ws.Range("H133").Value = 200000
$ws.Range("J133").Value = 200000
$ws.Range("L133").Value = 200000
$ws.Range("N133").Value = -205060

$ws.Range("H136").Value = 3926.8
$ws.Range("I136").Value = 3343.3333
$ws.Range("K136").Value = 10029.9999
$ws.Range("M136").Value = -7479.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3771.6667
$ws.Range("I20").Value = 2749.3333
$ws.Range("K20").Value = 2749.3333
$ws.Range("M20").Value = -2502.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 473
$ws.Range("I22").Value = 473
$ws.Range("K22").Value = 473
$ws.Range("M22").Value = -123

$ws.Range("H31").Value = 3728.0576
$ws.Range("I31").Value = 1542.5
$ws.Range("K31").Value = 1542.5
$ws.Range("M31").Value = -1247.5

$ws.Range("H34").Value = 3728.0576
$ws.Range("I34").Value = 1542.5
$ws.Range("K34").Value = 1542.5
$ws.Range("M34").Value = -1340.5

$ws.Range("H58").Value = 2810.8823
$ws.Range("I58").Value = 2658.1365
$ws.Range("K58").Value = 2658.1365
$ws.Range("M58").Value = -2455.1365

$ws.Range("H86").Value = 17986.666
$ws.Range("I86").Value = 17203.5
$ws.Range("K86").Value = 17203.5
$ws.Range("M86").Value = -16080.5

$ws.Range("H89").Value = 17986.666
$ws.Range("I89").Value = 17203.5
$ws.Range("K89").Value = 86017.5
$ws.Range("M89").Value = -80401.5

$ws.Range("H92").Value = 73142.5
$ws.Range("J92").Value = 73142.5
$ws.Range("L92").Value = 73142.5
$ws.Range("N92").Value = -78134.5

$ws.Range("H122").Value = 5969.0713
$ws.Range("I122").Value = 5780.3335
$ws.Range("J122").Value = 6110.625
$ws.Range("K122").Value = 17341.0005
$ws.Range("L122").Value = 18331.875
$ws.Range("M122").Value = -14891.0005
$ws.Range("N122").Value = -23231.875

$ws.Range("H132").Value = 3734.4666
$ws.Range("I132").Value = 3166.9092
$ws.Range("J132").Value = 5295.25
$ws.Range("K132").Value = 9500.7276
$ws.Range("L132").Value = 15885.75
$ws.Range("M132").Value = -6970.7276
$ws.Range("N132").Value = -20945.75

$ws.Range("H134").Value = 2886.111
$ws.Range("I134").Value = 2075.2
$ws.Range("J134").Value = 3899.75
$ws.Range("K134").Value = 6225.599999999999
$ws.Range("L134").Value = 11699.25
$ws.Range("M134").Value = -3690.599999999999
$ws.Range("N134").Value = -16769.25

$ws.Range("H136").Value = 2810.8823
$ws.Range("I136").Value = 2658.1365
$ws.Range("K136").Value = 7974.4095
$ws.Range("M136").Value = -5424.4095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 841.7778
$ws.Range("I68").Value = 696
$ws.Range("K68").Value = 2088
$ws.Range("M68").Value = -1277

$ws.Range("H71").Value = 841.7778
$ws.Range("I71").Value = 696
$ws.Range("K71").Value = 6264
$ws.Range("M71").Value = -2208

$ws.Range("H113").Value = 4049
$ws.Range("J113").Value = 4110.5557
$ws.Range("L113").Value = 12331.6671
$ws.Range("N113").Value = -16671.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 29888
$ws.Range("J93").Value = 29888
$ws.Range("L93").Value = 29888
$ws.Range("N93").Value = -33632

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 3557.5483
$ws.Range("I132").Value = 3051.7144
$ws.Range("K132").Value = 9155.143199999999
$ws.Range("M132").Value = -6625.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1002
$ws.Range("I10").Value = 503
$ws.Range("K10").Value = 503
$ws.Range("M10").Value = -363

$ws.Range("H14").Value = 11004
$ws.Range("I14").Value = 11004
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 11004
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -10832
$ws.Range("N14").ClearContents()

$ws.Range("H22").Value = 3455.8572
$ws.Range("I22").Value = 1617
$ws.Range("J22").Value = 4835
$ws.Range("K22").Value = 1617
$ws.Range("L22").Value = 4835
$ws.Range("M22").Value = -1322
$ws.Range("N22").Value = -5425

$ws.Range("H27").Value = 3455.8572
$ws.Range("I27").Value = 1617
$ws.Range("J27").Value = 4835
$ws.Range("K27").Value = 1617
$ws.Range("L27").Value = 4835
$ws.Range("M27").Value = -1510
$ws.Range("N27").Value = -5049

$ws.Range("H40").Value = 2548
$ws.Range("I40").Value = 2461.8
$ws.Range("J40").Value = 2979
$ws.Range("K40").Value = 2461.8
$ws.Range("L40").Value = 2979
$ws.Range("M40").Value = -2325.8
$ws.Range("N40").Value = -3251

$ws.Range("H46").Value = 7012.8276
$ws.Range("J46").Value = 7195.4287
$ws.Range("L46").Value = 7195.4287
$ws.Range("N46").Value = -7571.4287

$ws.Range("H55").Value = 3305.5
$ws.Range("I55").Value = 3241.5
$ws.Range("J55").Value = 3326.8333
$ws.Range("K55").Value = 3241.5
$ws.Range("L55").Value = 3326.8333
$ws.Range("M55").Value = -3068.5
$ws.Range("N55").Value = -3672.8333

$ws.Range("H61").Value = 1203.6666
$ws.Range("I61").Value = 1209.4
$ws.Range("K61").Value = 1209.4
$ws.Range("M61").Value = -1007.4

$ws.Range("H113").Value = 1203.6666
$ws.Range("I113").Value = 1209.4
$ws.Range("K113").Value = 1209.4
$ws.Range("M113").Value = 960.5999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 83322.336
$ws.Range("I10").Value = 50000
$ws.Range("K10").Value = 50000
$ws.Range("M10").Value = -49831

$ws.Range("H41").Value = 200005980
$ws.Range("J41").Value = 7472
$ws.Range("L41").Value = 7472
$ws.Range("N41").Value = -8252

$ws.Range("H81").Value = 4633.1113
$ws.Range("I81").Value = 3474.75
$ws.Range("J81").Value = 5559.8
$ws.Range("K81").Value = 6949.5
$ws.Range("L81").Value = 11119.6
$ws.Range("M81").Value = -5888.5
$ws.Range("N81").Value = -13241.6

$ws.Range("H84").Value = 4633.1113
$ws.Range("I84").Value = 3474.75
$ws.Range("J84").Value = 5559.8
$ws.Range("K84").Value = 34747.5
$ws.Range("L84").Value = 55598
$ws.Range("M84").Value = -29443.5
$ws.Range("N84").Value = -66206

$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

$ws.Range("H126").Value = 4667.6665
$ws.Range("I126").Value = 4004
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 12012
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -9542
$ws.Range("N126").Value = -19938.5

$ws.Range("H136").Value = 1544.6
$ws.Range("I136").Value = 1555.2858
$ws.Range("K136").Value = 4665.857400000001
$ws.Range("M136").Value = -2115.857400000001
